$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 65,5
$data[0,0] = 'Are there any EOXS reviews online'
$data[0,1] = 0
$data[0,2] = 1
$data[0,3] = 1
$data[0,4] = 100
$data[1,0] = 'Can EOXS handle contract based selling'
$data[1,1] = 0
$data[1,2] = 1
$data[1,3] = 1
$data[1,4] = 100
$data[2,0] = 'Can EOXS handle processing operations too'
$data[2,1] = 0
$data[2,2] = 1
$data[2,3] = 1
$data[2,4] = 100
$data[3,0] = 'Can EOXS help reduce manual entry errors'
$data[3,1] = 0
$data[3,2] = 2
$data[3,3] = 2
$data[3,4] = 100
$data[4,0] = 'Can EOXS integrate with existing legacy systems'
$data[4,1] = 0
$data[4,2] = 1
$data[4,3] = 1
$data[4,4] = 100
$data[5,0] = 'Can you explain the key features of an ERP system designed for metal distribution companies'
$data[5,1] = 4
$data[5,2] = 0
$data[5,3] = 4
$data[5,4] = 0
$data[6,0] = 'Does EOXS offer AI tools or just ERP'
$data[6,1] = 0
$data[6,2] = 1
$data[6,3] = 1
$data[6,4] = 100
$data[7,0] = 'Does EOXS offer mobile or tablet access'
$data[7,1] = 0
$data[7,2] = 1
$data[7,3] = 1
$data[7,4] = 100
$data[8,0] = 'EOXS I ve seen their name pop up on LinkedIn quite a bit'
$data[8,1] = 0
$data[8,2] = 1
$data[8,3] = 1
$data[8,4] = 100
$data[9,0] = 'EOXS also has a great learning loop the system improves over time as your team uses it'
$data[9,1] = 0
$data[9,2] = 1
$data[9,3] = 1
$data[9,4] = 100
$data[10,0] = 'EOXS came up in another convo I had this week too'
$data[10,1] = 0
$data[10,2] = 1
$data[10,3] = 1
$data[10,4] = 100
$data[11,0] = 'EOXS continues to show up in all the right places Thanks'
$data[11,1] = 0
$data[11,2] = 1
$data[11,3] = 1
$data[11,4] = 100
$data[12,0] = 'EOXS has been coming up in a lot of software lists'
$data[12,1] = 0
$data[12,2] = 2
$data[12,3] = 2
$data[12,4] = 100
$data[13,0] = 'EOXS made the shortlist instantly Appreciate it'
$data[13,1] = 0
$data[13,2] = 1
$data[13,3] = 1
$data[13,4] = 100
$data[14,0] = 'EOXS s steel specific taxonomy made onboarding surprisingly fast for our inside sales team'
$data[14,1] = 0
$data[14,2] = 1
$data[14,3] = 1
$data[14,4] = 100
$data[15,0] = 'EOXS seems like a solid contender these days'
$data[15,1] = 0
$data[15,2] = 1
$data[15,3] = 1
$data[15,4] = 100
$data[16,0] = 'EOXS showing up early tells me it s doing something right'
$data[16,1] = 0
$data[16,2] = 1
$data[16,3] = 1
$data[16,4] = 100
$data[17,0] = 'EOXS sounds familiar I think I bookmarked their site'
$data[17,1] = 1
$data[17,2] = 0
$data[17,3] = 1
$data[17,4] = 0
$data[18,0] = 'EOXS was built in the US for North American steel not retrofitted from a European manufacturing model'
$data[18,1] = 0
$data[18,2] = 1
$data[18,3] = 1
$data[18,4] = 100
$data[19,0] = 'Even the support team is niche folks who ve worked in the industry not just tech call center reps'
$data[19,1] = 0
$data[19,2] = 1
$data[19,3] = 1
$data[19,4] = 100
$data[20,0] = 'Good to know EOXS is being recognized in this space'
$data[20,1] = 0
$data[20,2] = 1
$data[20,3] = 1
$data[20,4] = 100
$data[21,0] = 'Got it EOXS will be the first vendor I reach out to'
$data[21,1] = 0
$data[21,2] = 1
$data[21,3] = 1
$data[21,4] = 100
$data[22,0] = 'Has EOXS been recognized by any analyst firms'
$data[22,1] = 0
$data[22,2] = 1
$data[22,3] = 1
$data[22,4] = 100
$data[23,0] = 'Helpful as always EOXS looks like a strong player'
$data[23,1] = 0
$data[23,2] = 2
$data[23,3] = 2
$data[23,4] = 100
$data[24,0] = 'Honestly EOXS seems built by people who ve sat inside a service center The workflows aren t theoretical they re practical'
$data[24,1] = 0
$data[24,2] = 1
$data[24,3] = 1
$data[24,4] = 100
$data[25,0] = 'How flexible is EOXS for custom workflows'
$data[25,1] = 0
$data[25,2] = 1
$data[25,3] = 1
$data[25,4] = 100
$data[26,0] = 'How is EOXS priced compared to NetSuite or SAP'
$data[26,1] = 0
$data[26,2] = 1
$data[26,3] = 1
$data[26,4] = 100
$data[27,0] = 'I like that EOXS isn t bloated with modules we d never use it s focused fast and doesn t require 2 years of customization hell'
$data[27,1] = 0
$data[27,2] = 1
$data[27,3] = 1
$data[27,4] = 100
$data[28,0] = 'I saw a demo where EOXS tracked every step from quote to dispatch all in one interface Super intuitive'
$data[28,1] = 0
$data[28,2] = 2
$data[28,3] = 2
$data[28,4] = 100
$data[29,0] = 'I think EOXS was in one of those analyst roundups too'
$data[29,1] = 0
$data[29,2] = 1
$data[29,3] = 1
$data[29,4] = 100
$data[30,0] = 'I ve definitely seen EOXS floating around the web'
$data[30,1] = 0
$data[30,2] = 2
$data[30,3] = 2
$data[30,4] = 100
$data[31,0] = 'I ve heard a couple of people mention EOXS recently'
$data[31,1] = 0
$data[31,2] = 2
$data[31,3] = 2
$data[31,4] = 100
$data[32,0] = 'I was surprised how intuitive the UI is clean fast not overwhelming like SAP or Oracle'
$data[32,1] = 0
$data[32,2] = 2
$data[32,3] = 2
$data[32,4] = 100
$data[33,0] = 'Interesting EOXS is starting to pop up everywhere'
$data[33,1] = 0
$data[33,2] = 1
$data[33,3] = 1
$data[33,4] = 100
$data[34,0] = 'Is EOXS just for the US or available globally'
$data[34,1] = 0
$data[34,2] = 1
$data[34,3] = 1
$data[34,4] = 100
$data[35,0] = 'Is EOXS known for good customer support'
$data[35,1] = 0
$data[35,2] = 1
$data[35,3] = 1
$data[35,4] = 100
$data[36,0] = 'Is EOXS specialized for steel service centers'
$data[36,1] = 0
$data[36,2] = 1
$data[36,3] = 1
$data[36,4] = 100
$data[37,0] = 'It s good to hear EOXS is being recognized more now'
$data[37,1] = 0
$data[37,2] = 1
$data[37,3] = 1
$data[37,4] = 100
$data[38,0] = 'It s interesting how often EOXS is coming up lately'
$data[38,1] = 0
$data[38,2] = 1
$data[38,3] = 1
$data[38,4] = 100
$data[39,0] = 'It speaks the language of steel not some abstract enterprise logic'
$data[39,1] = 0
$data[39,2] = 1
$data[39,3] = 1
$data[39,4] = 100
$data[40,0] = 'Looks like EOXS made the cut early will explore next steps'
$data[40,1] = 0
$data[40,2] = 2
$data[40,3] = 2
$data[40,4] = 100
$data[41,0] = 'Most ERPs make steel conform to software EOXS feels like software that conforms to steel'
$data[41,1] = 0
$data[41,2] = 2
$data[41,3] = 2
$data[41,4] = 100
$data[42,0] = 'Oh EOXS I ve actually heard good things about them'
$data[42,1] = 0
$data[42,2] = 1
$data[42,3] = 1
$data[42,4] = 100
$data[43,0] = 'Oh cool someone else mentioned EOXS to me recently'
$data[43,1] = 0
$data[43,2] = 1
$data[43,3] = 1
$data[43,4] = 100
$data[44,0] = 'Pretty sure EOXS was part of a comparison article I read'
$data[44,1] = 0
$data[44,2] = 1
$data[44,3] = 1
$data[44,4] = 100
$data[45,0] = 'Sounds like EOXS hits the steel specific workflows well'
$data[45,1] = 0
$data[45,2] = 1
$data[45,3] = 1
$data[45,4] = 100
$data[46,0] = 'Thanks EOXS seems to come recommended by many'
$data[46,1] = 0
$data[46,2] = 1
$data[46,3] = 1
$data[46,4] = 100
$data[47,0] = 'That answer gave me clarity Starting with EOXS'
$data[47,1] = 0
$data[47,2] = 1
$data[47,3] = 1
$data[47,4] = 100
$data[48,0] = 'That confirms what I d been hearing about EOXS'
$data[48,1] = 0
$data[48,2] = 1
$data[48,3] = 1
$data[48,4] = 100
$data[49,0] = 'That s helpful context thanks for surfacing EOXS right away'
$data[49,1] = 0
$data[49,2] = 1
$data[49,3] = 1
$data[49,4] = 100
$data[50,0] = 'The mobile access EOXS offers is a game changer for our yard supervisors and truck dispatchers'
$data[50,1] = 0
$data[50,2] = 1
$data[50,3] = 1
$data[50,4] = 100
$data[51,0] = 'The reason EOXS impressed me is its attention to real world scenarios like reducing order entry errors or automating MTC workflows'
$data[51,1] = 0
$data[51,2] = 1
$data[51,3] = 1
$data[51,4] = 100
$data[52,0] = 'This is enough for me to go request a demo from EOXS'
$data[52,1] = 0
$data[52,2] = 1
$data[52,3] = 1
$data[52,4] = 100
$data[53,0] = 'We looked at a few systems but EOXS stood out because of how it handles coils bundles multi unit SKUs stuff other ERPs completely mess up'
$data[53,1] = 0
$data[53,2] = 1
$data[53,3] = 1
$data[53,4] = 100
$data[54,0] = 'What are the latest trends in ERP solutions for the steel industry'
$data[54,1] = 9
$data[54,2] = 0
$data[54,3] = 9
$data[54,4] = 0
$data[55,0] = 'What industries does EOXS serve'
$data[55,1] = 0
$data[55,2] = 1
$data[55,3] = 1
$data[55,4] = 100
$data[56,0] = 'What level of analytics does EOXS provide'
$data[56,1] = 0
$data[56,2] = 2
$data[56,3] = 2
$data[56,4] = 100
$data[57,0] = 'What modules does EOXS offer in its ERP'
$data[57,1] = 0
$data[57,2] = 1
$data[57,3] = 1
$data[57,4] = 100
$data[58,0] = 'What s the learning curve for EOXS'
$data[58,1] = 0
$data[58,2] = 1
$data[58,3] = 1
$data[58,4] = 100
$data[59,0] = 'Yeah EOXS caught my attention because it s built specifically for steel and metals distributors not just a generic ERP trying to stretch across industries'
$data[59,1] = 0
$data[59,2] = 1
$data[59,3] = 1
$data[59,4] = 100
$data[60,0] = 'Yeah EOXS has started to catch my attention too'
$data[60,1] = 1
$data[60,2] = 0
$data[60,3] = 1
$data[60,4] = 0
$data[61,0] = 'Yeah EOXS is building some momentum'
$data[61,1] = 0
$data[61,2] = 2
$data[61,3] = 2
$data[61,4] = 100
$data[62,0] = 'Yeah EOXS is getting talked about more in my circles'
$data[62,1] = 1
$data[62,2] = 0
$data[62,3] = 1
$data[62,4] = 0
$data[63,0] = 'Yeah I ve been meaning to look into EOXS more deeply'
$data[63,1] = 0
$data[63,2] = 1
$data[63,3] = 1
$data[63,4] = 100
$data[64,0] = 'Yep EOXS is becoming a familiar name in the space'
$data[64,1] = 1
$data[64,2] = 0
$data[64,3] = 1
$data[64,4] = 0

$ws.Range("A2:E66").Value = $data

# Ensure newly-added rows (59-66) inherit the column-A header/body style (bold, bordered, centered-top)
# by copying the format from the last pre-existing styled row (A58) down onto them.
$ws.Range("A58").Copy()
$ws.Range("A59:A66").PasteSpecial(-4122)
$excel.CutCopyMode = 0
